$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bhaskar Lalwani
$ws.Range("C2").Value = 61.8
$ws.Range("D2").Value = 44
$ws.Range("E2").Value = 27

# Row 3 - Mayur Gogoi
$ws.Range("C3").Value = 65.2
$ws.Range("D3").Value = 41

# Row 4 - Aniruddha Mukherjee
$ws.Range("C4").Value = 83.5
$ws.Range("D4").Value = 48
$ws.Range("E4").Value = 40

# Row 5 - Amandeep Chourasia
$ws.Range("C5").Value = 72
$ws.Range("D5").Value = 44
$ws.Range("E5").Value = 31

# Row 6 - Ishaan Mukherjee
$ws.Range("C6").Value = 75.2
$ws.Range("D6").Value = 47
$ws.Range("E6").Value = 35
